# ============================================================
# Edit: add "2022-Q4" sheet (with fund holdings data) as the
# second sheet, and prepend a matching summary row on "总计".
# ============================================================

$wb = $excel.ActiveWorkbook

# ---- 1. Insert the new "2022-Q4" worksheet right before "2022-Q3" ----
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($sheetQ3)
$q4.Name = "2022-Q4"

# Columns B:G hold text-like values (fund codes, names, and numbers
# that are stored as text in the source data) - force text format
# up front so values such as "159941" or "000043" keep their shape
# (leading zeros) instead of being coerced into numbers.
$q4.Range("B1:G26").NumberFormat = "@"

# Header row
$q4Header = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($j = 0; $j -lt $q4Header.Length; $j++) {
    $q4.Cells.Item(1, 2 + $j).Value = $q4Header[$j]
}

# Data rows (fund code, fund name, fund size, total equity position,
# position ratio, held market value, position rank)
$q4Data = @(
    @("159941", "广发纳斯达克100ETF（QDII）", "114.77", "90.42", "1.99", "2.2839", "10"),
    @("513100", "国泰纳斯达克100（QDII-ETF）", "51.50", "90.72", "2.08", "1.0712", "9"),
    @("159632", "华安纳斯达克100ETF（QDII）", "30.39", "87.77", "1.82", "0.5531", "10"),
    @("160213", "国泰纳斯达克100指数（QDII）", "15.65", "90.80", "1.89", "0.2958", "10"),
    @("000834", "大成纳斯达克100指数（QDII）", "15.61", "81.77", "1.69", "0.2638", "10"),
    @("513300", "华夏纳斯达克100ETF（QDII）", "12.43", "97.54", "2.03", "0.2523", "8"),
    @("000043", "嘉实美国成长股票（QDII）人民币", "12.69", "92.23", "1.46", "0.1853", "10"),
    @("000044", "嘉实美国成长股票（QDII）美元现汇", "12.69", "92.23", "1.46", "0.1853", "10"),
    @("161130", "易方达纳斯达克100指数人民币（QDII-LOF）", "7.77", "90.34", "1.88", "0.1461", "10"),
    @("003722", "易方达纳斯达克100指数美元（QDII-LOF）A", "7.77", "90.34", "1.88", "0.1461", "10"),
    @("012868", "易方达标普信息科技指数（QDII-LOF）人民币 C", "5.09", "91.36", "2.50", "0.1272", "6"),
    @("161128", "易方达标普信息科技指数（QDII-LOF）人民币", "5.09", "91.36", "2.50", "0.1272", "6"),
    @("003721", "易方达标普信息科技指数（QDII-LOF）美元A", "4.93", "91.36", "2.50", "0.1232", "6"),
    @("016532", "嘉实纳斯达克100指数（QDII）A人民币", "1.12", "94.67", "1.97", "0.0221", "10"),
    @("016533", "嘉实纳斯达克100指数（QDII）C人民币", "1.12", "94.67", "1.97", "0.0221", "10"),
    @("016534", "嘉实纳斯达克100指数（QDII）A美元现汇", "1.12", "94.67", "1.97", "0.0221", "10"),
    @("016535", "嘉实纳斯达克100指数（QDII）C美元现汇", "1.12", "94.67", "1.97", "0.0221", "10"),
    @("016055", "博时纳斯达克100指数（QDII）A人民币", "1.06", "90.62", "1.88", "0.0199", "10"),
    @("016057", "博时纳斯达克100指数（QDII）C人民币", "1.06", "90.62", "1.88", "0.0199", "10"),
    @("016056", "博时纳斯达克100指数（QDII）A美元现汇", "1.06", "90.62", "1.88", "0.0199", "10"),
    @("016058", "博时纳斯达克100指数（QDII）C美元现汇", "1.06", "90.62", "1.88", "0.0199", "10"),
    @("005698", "华夏全球科技先锋混合（QDII）", "0.60", "83.35", "2.95", "0.0177", "10"),
    @("012869", "易方达标普信息科技指数（QDII-LOF）美元 C", "0.16", "91.36", "2.50", "0.0040", "6"),
    @("012870", "易方达纳斯达克100指数人民币（QDII-LOF）C", "0.21", "90.34", "1.88", "0.0039", "10"),
    @("012871", "易方达纳斯达克100指数美元（QDII-LOF）C", "0.21", "90.34", "1.88", "0.0039", "10")
)

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $row = $q4Data[$i]
    $r = 2 + $i
    $q4.Cells.Item($r, 1).Value = $i
    for ($j = 0; $j -lt 6; $j++) {
        $q4.Cells.Item($r, 2 + $j).Value = $row[$j]
    }
    $q4.Cells.Item($r, 8).Value = [int]$row[6]
}

# Header style: bold, centered, thin-bordered (matches the look of
# the header rows on the other quarterly sheets)
$q4HeaderRange = $q4.Range("B1:H1")
$q4HeaderRange.Font.Bold = $true
$q4HeaderRange.HorizontalAlignment = -4108
$q4HeaderRange.VerticalAlignment = -4160
$q4HeaderRange.Borders.LineStyle = 1
$q4HeaderRange.Borders.Weight = 2

# Row-index column (A) style: same bold/centered/bordered look
$q4IndexRange = $q4.Range("A2:A26")
$q4IndexRange.Font.Bold = $true
$q4IndexRange.HorizontalAlignment = -4108
$q4IndexRange.VerticalAlignment = -4160
$q4IndexRange.Borders.LineStyle = 1
$q4IndexRange.Borders.Weight = 2

# ---- 2. Prepend the 2022-Q4 summary row on the "总计" sheet ----
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()
$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 25
$summary.Range("D2").Value = 5.96

$summaryA2 = $summary.Range("A2")
$summaryA2.Font.Bold = $true
$summaryA2.HorizontalAlignment = -4108
$summaryA2.VerticalAlignment = -4160
$summaryA2.Borders.LineStyle = 1
$summaryA2.Borders.Weight = 2

# Renumber the index column for the rows that got pushed down
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

Write-Host "Edit applied: 2022-Q4 sheet added and 总计 summary updated"
